# LinkedInUserData.xlsx - swap in a new batch of page-factory test accounts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Islam Abd Alazez -> Waleed Ali (islamtalkha83@gmail.com)
$ws.Range("A2").Value = "Waleed"
$ws.Range("B3").Value = "fadel"
$ws.Range("B4").Value = "Gamal"
$ws.Range("C3").Value = "ali89afit@gmail.com"
$ws.Range("C4").Value = "gamal79afit@gmail.com"
$ws.Range("C2").Value = "islamtalkha83@gmail.com"
$ws.Range("B2").Value = "Ali"
$ws.Range("A3").Value = "Ali"
$ws.Range("A4").Value = "Ahmed"

# Widen the E-mail column to fit the new, longer addresses.
$ws.Columns("C").ColumnWidth = 25.71

# Move the active selection to C2 (where the new first e-mail now lives).
$ws.Range("C2").Select()
